# Weekly update: insert a new price record as the first data row (row 28)
# for "Vega Monumental Concepción - Poroto verde", pushing the existing
# rows 28-39 down to 29-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28, shifting rows 28:39 down to 29:40.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with this week's record.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44609
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112031
$ws.Cells.Item(28, 7).Value = "Poroto verde"
$ws.Cells.Item(28, 8).Value = "Magnum"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 28000
$ws.Cells.Item(28, 12).Value = 30000
$ws.Cells.Item(28, 13).Value = 29000
$ws.Cells.Item(28, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 1160
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"

# Match the date style used by the rest of the "Fecha" column.
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
